$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = -7.912099999999997
$ws.Range("D7").Value = -7.344299999999995
$ws.Range("B8").Value = 5.040600000000005
$ws.Range("A12").Value = -22.86590000000001
$ws.Range("B12").Value = 5.328
$ws.Range("B14").Value = 8.809500000000007
$ws.Range("D19").Value = -8.602799999999993
$ws.Range("D21").Value = -7.540599999999996
$ws.Range("B22").Value = 4.842400000000005
$ws.Range("D24").Value = -8.06819999999999
